$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting existing rows 8-97 down to 9-98.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new price record.
# Columns A, B, C, E, F, G, I, Q, R are constant across the whole table.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44817
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112031
$ws.Cells.Item(8, 7).Value = "Poroto verde"
$ws.Cells.Item(8, 8).Value = "Magnum"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 30000
$ws.Cells.Item(8, 12).Value = 32000
$ws.Cells.Item(8, 13).Value = 31000
$ws.Cells.Item(8, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 1240
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
